# Apply the "added a fixed bug comment" edit to the test_android.xlsx bug
# tracking sheet.
#
# Semantic change: cell C5 ("还未能重现") gets a longer follow-up comment
# appended describing the root cause / fix. Everything else in the diff
# (shared-string reindexing, font "family" normalisation, row 5 growing
# taller) is a mechanical side effect of that one content edit plus a
# couple of cosmetic formatting touches (wrap-text on C1/C4/C5, the
# selection left on A4:A5, and turning on portrait page setup).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C5: replace the short "still can't reproduce" note with the fuller
#     follow-up comment that explains the root cause and the fix. ---
$ws.Cells.Item(5, 3).Value = "还未能重现。但是umeng那边抓到了错误信息，是""NullPointerException: replacement == null""。已经把这一处修改了。"

# --- wrap text formatting on the column-3 notes cells that should carry
#     style 4 (wrap text) instead of style 1 (plain) ---
$ws.Cells.Item(1, 3).WrapText = $true
$ws.Cells.Item(4, 3).WrapText = $true
$ws.Cells.Item(5, 3).WrapText = $true

# Row 5 grows to fit the longer wrapped comment.
$ws.Rows.Item(5).RowHeight = 52.2

# Selection moves to the merged Bug cell (A4:A5) that anchors the row just
# edited.
[void]$ws.Range("A4:A5").Select()

# Turn on (portrait) page setup for printing.
$ws.PageSetup.Orientation = 1
